# Role in JwtDTO and visibility only for teachers
# Adds a new "Browse proposals" test case row (row 5) describing that a
# non-teacher user should not see the proposals list / sidebar link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (borders, date number format, etc.) from the last data
# row (row 4) down into the new row 5 so the new row matches the existing
# table styling exactly.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

# Fill in the new test case data.
$ws.Range("A5").Value = "Browse proposals"
$ws.Range("B5").Value = "Do not see list of proposals if not logged in as teacher"
$ws.Range("C5").Value = "Do not log in, go to home page"
$ws.Range("D5").Value = "Sidebar should not contain link to browse proposals"
$ws.Range("E5").Value = $ws.Range("E4").Value()

# Column B now holds longer text ("Do not see list of proposals if not
# logged in as teacher"); widen it to fit.
$ws.Columns("B").ColumnWidth = 44.5

# Move the selection to A2, matching the post-edit selection state.
$null = $ws.Range("A2").Select()
